$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "A" column (distance-without-correction raw values);
# its data is being rebuilt in column C below.
$ws.Range("A1:A21").Clear()

# Header row: new "Original distance" / moon-distance columns added,
# existing headers shifted right by one column.
$ws.Range("B1").Value = "Original distance"
$ws.Range("C1").Value = "Distance without correction"
$ws.Range("D1").Value = "Distance with temperature and huminity correction"
$ws.Range("E1").Value = "moon Distance"
$ws.Range("F1").Value = "moon distance with correction"

# Data rows
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 3.7

$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 4.8
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 4.5

$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 6.2
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 6.31

$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 8.6
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 7.91

$ws.Range("B6").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 10.1
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 9.38

$ws.Range("B7").Value = 12
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 12.6
$ws.Range("E7").Value = 11
$ws.Range("F7").Value = 11.39

$ws.Range("B8").Value = 14
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 14.5
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 13.82

$ws.Range("B9").Value = 16
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 16.5
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 15.56

$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 18
$ws.Range("D10").Value = 18.5
$ws.Range("E10").Value = 18
$ws.Range("F10").Value = 18.2

$ws.Range("B11").Value = 20
$ws.Range("C11").Value = 20
$ws.Range("D11").Value = 20.3
$ws.Range("E11").Value = 21
$ws.Range("F11").Value = 21.05

$ws.Range("B12").Value = 22
$ws.Range("C12").Value = 22
$ws.Range("D12").Value = 22.6
$ws.Range("E12").Value = 22
$ws.Range("F12").Value = 22.09

$ws.Range("B13").Value = 24
$ws.Range("C13").Value = 24
$ws.Range("D13").Value = 24.1
$ws.Range("E13").Value = 24
$ws.Range("F13").Value = 23.97

$ws.Range("B14").Value = 26
$ws.Range("C14").Value = 26
$ws.Range("D14").Value = 26.6
$ws.Range("E14").Value = 25
$ws.Range("F14").Value = 25.67

$ws.Range("B15").Value = 28
$ws.Range("C15").Value = 28
$ws.Range("D15").Value = 28.6
$ws.Range("E15").Value = 27
$ws.Range("F15").Value = 27.45

$ws.Range("B16").Value = 30
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 29.8
$ws.Range("E16").Value = 30
$ws.Range("F16").Value = 30.02

$ws.Range("B17").Value = 32
$ws.Range("C17").Value = 32
$ws.Range("D17").Value = 31.7
$ws.Range("E17").Value = 30
$ws.Range("F17").Value = 30.09

$ws.Range("B18").Value = 34
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = 34.4
$ws.Range("E18").Value = 33
$ws.Range("F18").Value = 33.36

$ws.Range("B19").Value = 36
$ws.Range("C19").Value = 36
$ws.Range("D19").Value = 36.3
$ws.Range("E19").Value = 33
$ws.Range("F19").Value = 33.64

$ws.Range("B20").Value = 38
$ws.Range("C20").Value = 38
$ws.Range("D20").Value = 37.8
$ws.Range("E20").Value = 34
$ws.Range("F20").Value = 35.65

$ws.Range("B21").Value = 40
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 39.6
$ws.Range("E21").Value = 36
$ws.Range("F21").Value = 36.79

# Restore the active-cell selection recorded in the workbook view
$ws.Range("G18").Select() | Out-Null
